$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2/B3 values (step 2.2 intraday data labels)
$ws.Range("B2").Value = "v2.2"
$ws.Range("B3").Value = "yup"

# Append new row 28 of data
$ws.Range("B28").Value = "lastrow"
$ws.Range("C28").Value = "-"
$ws.Range("D28").Value = -3.3
$ws.Range("E28").Formula = "=A28*D28"
$ws.Range("F28").Value = "x"

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:F28"))

# Update the selection to reflect where the user ended up
$ws.Range("D29").Select() | Out-Null
